# Add the missing project row (row 3) so that all projects are shown.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values look numeric need to be forced to text format first,
# otherwise Excel will auto-convert the typed string into a real number.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("P3").NumberFormat = "@"

$ws.Range("A3").Value = "1744625747268"
$ws.Range("B3").Value = "Caterpillar"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "Hybrid Seat Allocation"
$ws.Range("E3").Value = "PI"
$ws.Range("F3").Value = "Copi"
$ws.Range("G3").Value = "2023-2024"
$ws.Range("H3").Value = "20000"
$ws.Range("I3").Value = "10000"
$ws.Range("J3").Value = "Details"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "h"
$ws.Range("M3").Value = "agreementDocument-1744625747215-658333568.pdf"
$ws.Range("O3").Value = "2025-04-14T10:15:47.268Z"
$ws.Range("P3").Value = "1744624259342"
